$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Start Game Selected" -> three separate paragraphs: "Selects:", "Start",
#    "Exit"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Start Game Selected") | Out-Null
$r.InsertXML("<w:p><w:r><w:t>Selects:</w:t></w:r></w:p><w:p><w:r><w:t>Start</w:t></w:r></w:p><w:p><w:r><w:t>Exit</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 2) "Level" cell gains a second paragraph "Desktop"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Level") | Out-Null
$r.InsertXML("<w:p><w:r><w:t>Level</w:t></w:r></w:p><w:p><w:r><w:t>Desktop</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3) "Character move to direction that's pressed" cell gains a trailing
#    empty paragraph
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Character move to direction that" + [char]0x2019 + "s pressed") | Out-Null
$r.InsertXML("<w:p><w:r><w:t xml:space='preserve'>Character </w:t></w:r><w:r><w:t>move to direction that" + [char]0x2019 + "s pressed</w:t></w:r></w:p><w:p/>")

# ---------------------------------------------------------------------------
# 4) Merge the two runs of "After Jump (Keyboard Space Bar 360 Pad RT" / ")
#    Light Kick or Heavy Punch ( Keyboard (K or L 360 Pad)" into one run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("After Jump (Keyboard Space Bar 360 Pad RT) Light Kick or Heavy Punch ( Keyboard (K or L 360 Pad)", $true, $false, $false, $false, $false, $true, 1, $false, "After Jump (Keyboard Space Bar 360 Pad RT) Light Kick or Heavy Punch ( Keyboard (K or L 360 Pad)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "Character in dive and kick the enemy ..." -> remove " and" (keep the
#    two original runs separate)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Character in dive and kick the enemy") | Out-Null
$p = $r.Paragraphs(1)
$pr = $p.Range
$pr.InsertXML("<w:p><w:r><w:t xml:space='preserve'>Character in dive kick the enemy to give spacing between </w:t></w:r><w:r><w:t>enemy and character</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 6) "Running punch" -> "Cork Screw" + bookmark (_GoBack) + " punch"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Running punch") | Out-Null
$r.InsertXML("<w:p><w:r><w:t>Cork Screw</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/><w:r><w:t xml:space='preserve'> punch</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 7) "While Character crouch he can roll left or right for invincibility
#    frames" is split into two runs (second one starts with a
#    lastRenderedPageBreak) and the following "Game" cell also gains a
#    lastRenderedPageBreak before its text
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("While Character crouch he can roll left or right for invincibility frames") | Out-Null
$p = $r.Paragraphs(1)
$pr = $p.Range
$pr.InsertXML("<w:p><w:r><w:t xml:space='preserve'>While Character crouch he can roll left or right for </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>invincibility frames</w:t></w:r></w:p>")

$r2 = $d.Range($pr.End, $d.Content.End)
$r2.Find.Execute("Game") | Out-Null
$r2.InsertXML("<w:p><w:r><w:lastRenderedPageBreak/><w:t>Game</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 8) Merge the trailing three runs after "Light Kick(Keyboard K 360 Pad A)"
#    into a single run, keeping the first run separate
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Light Kick(Keyboard K 360 Pad A) 2x, Jump(Keyboard Spacebar 360 Pad RT), Heavy Kick (Keyboard L, 360 Pad B)") | Out-Null
$r.InsertXML("<w:p><w:r><w:t>Light Kick(Keyboard K 360 Pad A)</w:t></w:r><w:r><w:t xml:space='preserve'> 2x, Jump(Keyboard Spacebar 360 Pad RT), Heavy Kick (Keyboard L, 360 Pad B)</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 9) Remove the stray bookmark (_GoBack) after the Crouch/Heavy Punch combo
#    text, keeping the three runs intact
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Crouch (Keyboard Ctrl, 360 Pad LT), Light Punch 2x, Heavy Punch( Keyboard I 360 Pad  B)") | Out-Null
$p = $r.Paragraphs(1)
$pr = $p.Range
$pr.InsertXML("<w:p><w:r><w:t>Crouch (Keyboard Ctrl, 360 Pad LT), Light Punch 2x, Heavy Punch</w:t></w:r><w:r><w:t>( Keyboard I 360 Pad  B</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 10) "When the Character moves to the next area a number of enemies will
#     spawn" -> redistribute the text across the existing four runs, moving
#     "number of " to the start of the run carrying lastRenderedPageBreak
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("When the Character moves to the next area a number of enemies will spawn") | Out-Null
$r.InsertXML("<w:p><w:r><w:t xml:space='preserve'>When the Character moves to the next area a </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>number of en</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>mies will spawn</w:t></w:r></w:p>")
